$p = $ppt.ActivePresentation
$s9 = $p.Slides.Item(9)
$s9.Delete()
Write-Output ("count after delete: " + $p.Slides.Count)
$ppt.Undo()
Write-Output ("count after undo: " + $p.Slides.Count)
